# Rename the column headers in row 1 so they carry the respective
# input-file-version suffix instead of the generic _old / _new markers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

# Columns A-J (1-10): "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = 1 + $i
    $ws.Cells.Item(1, $col).Value = $baseHeaders[$i] + "_FV2404"
}

# Column K (11) "diff" is unchanged.

# Columns L-U (12-21): "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseHeaders.Count; $i++) {
    $col = 12 + $i
    $ws.Cells.Item(1, $col).Value = $baseHeaders[$i] + "_FV2410"
}

# Turn the data range into a proper Excel Table ("Table1").
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
